$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 666.25
$ws.Range("I2").Value = 666.25
$ws.Range("K2").Value = 666.25
$ws.Range("M2").Value = -553.25

$ws.Range("H4").Value = 9000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""

$ws.Range("H9").Value = 191.06667
$ws.Range("I9").Value = 145.3
$ws.Range("J9").Value = 282.6
$ws.Range("K9").Value = 145.3
$ws.Range("L9").Value = 282.6
$ws.Range("M9").Value = 23.69999999999999
$ws.Range("N9").Value = -620.6

$ws.Range("H10").Value = 1835
$ws.Range("I10").Value = 1750
$ws.Range("J10").Value = 2005
$ws.Range("K10").Value = 1750
$ws.Range("L10").Value = 2005
$ws.Range("M10").Value = -1457
$ws.Range("N10").Value = -2591

$ws.Range("H32").Value = 10499.625
$ws.Range("J32").Value = 10499.625
$ws.Range("L32").Value = 10499.625
$ws.Range("N32").Value = -11151.625

$ws.Range("H43").Value = 4167.6665
$ws.Range("I43").Value = 4500
$ws.Range("J43").Value = 4001.5
$ws.Range("K43").Value = 4500
$ws.Range("L43").Value = 4001.5
$ws.Range("M43").Value = -4431
$ws.Range("N43").Value = -4139.5

$ws.Range("H76").Value = 2624.5
$ws.Range("I76").Value = 2624.5
$ws.Range("K76").Value = 2624.5
$ws.Range("M76").Value = -2309.5

$ws.Range("H79").Value = 2624.5
$ws.Range("I79").Value = 2624.5
$ws.Range("K79").Value = 2624.5
$ws.Range("M79").Value = -1532.5

$ws.Range("H113").Value = 10000
$ws.Range("J113").Value = 10000
$ws.Range("L113").Value = 10000
$ws.Range("N113").Value = -16508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 622.5
$ws.Range("I5").Value = 496.66666
$ws.Range("K5").Value = 496.66666
$ws.Range("M5").Value = -384.66666

$ws.Range("H134").Value = 78615
$ws.Range("J134").Value = 78615
$ws.Range("L134").Value = 78615
$ws.Range("N134").Value = -88755

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 622.5
$ws.Range("I4").Value = 496.66666
$ws.Range("K4").Value = 496.66666
$ws.Range("M4").Value = -381.66666

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""

$ws.Range("H80").Value = 295
$ws.Range("J80").Value = 285
$ws.Range("L80").Value = 285
$ws.Range("N80").Value = -2281

$ws.Range("H83").Value = 295
$ws.Range("J83").Value = 285
$ws.Range("L83").Value = 1425
$ws.Range("N83").Value = -11409

$ws.Range("H94").Value = 3184.75
$ws.Range("I94").Value = 2870.7
$ws.Range("K94").Value = 2870.7
$ws.Range("M94").Value = -2419.7

$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175

$ws.Range("H111").Value = 98997
$ws.Range("J111").Value = 98997
$ws.Range("L111").Value = 98997
$ws.Range("N111").Value = -107177

$ws.Range("H114").Value = 94684
$ws.Range("J114").Value = 94684
$ws.Range("L114").Value = 94684
$ws.Range("N114").Value = -103362

$ws.Range("H115").Value = 90339.5
$ws.Range("J115").Value = 90339.5
$ws.Range("L115").Value = 90339.5
$ws.Range("N115").Value = -93473.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1558
$ws.Range("I31").Value = 1460
$ws.Range("J31").Value = 1586
$ws.Range("K31").Value = 1460
$ws.Range("L31").Value = 1586
$ws.Range("M31").Value = -1165
$ws.Range("N31").Value = -2176

$ws.Range("H34").Value = 1558
$ws.Range("I34").Value = 1460
$ws.Range("J34").Value = 1586
$ws.Range("K34").Value = 1460
$ws.Range("L34").Value = 1586
$ws.Range("M34").Value = -1258
$ws.Range("N34").Value = -1990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 400363.06
$ws.Range("I4").Value = 357531.94
$ws.Range("K4").Value = 1072595.82
$ws.Range("M4").Value = -1072483.82

$ws.Range("H26").Value = 501
$ws.Range("I26").Value = 501
$ws.Range("K26").Value = 1503
$ws.Range("M26").Value = -1215

$ws.Range("H39").Value = 3587.75
$ws.Range("I39").Value = 2231.3333
$ws.Range("J39").Value = 4401.6
$ws.Range("K39").Value = 6693.999899999999
$ws.Range("L39").Value = 13204.8
$ws.Range("M39").Value = -6399.999899999999
$ws.Range("N39").Value = -13792.8

$ws.Range("H103").Value = 2842
$ws.Range("I103").Value = 3174.6667
$ws.Range("J103").Value = 2509.3333
$ws.Range("K103").Value = 9524.000100000001
$ws.Range("L103").Value = 7527.999899999999
$ws.Range("M103").Value = -8645.000100000001
$ws.Range("N103").Value = -9285.999899999999

$ws.Range("H129").Value = 2541.15
$ws.Range("I129").Value = 2124
$ws.Range("J129").Value = 2819.25
$ws.Range("K129").Value = 6372
$ws.Range("L129").Value = 8457.75
$ws.Range("M129").Value = -1372
$ws.Range("N129").Value = -18457.75

$ws.Range("H136").Value = 4625
$ws.Range("I136").Value = 4625
$ws.Range("K136").Value = 13875
$ws.Range("M136").Value = -8775

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 75.625
$ws.Range("I2").Value = 75.625
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 75.625
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 37.375
$ws.Range("N2").Value = ""

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 10003
$ws.Range("I14").Value = 7502
$ws.Range("K14").Value = 7502
$ws.Range("M14").Value = -7330

$ws.Range("H22").Value = 2683.611
$ws.Range("I22").Value = 1771.1428
$ws.Range("J22").Value = 3264.2727
$ws.Range("K22").Value = 1771.1428
$ws.Range("L22").Value = 3264.2727
$ws.Range("M22").Value = -1476.1428
$ws.Range("N22").Value = -3854.2727

$ws.Range("H27").Value = 2683.611
$ws.Range("I27").Value = 1771.1428
$ws.Range("J27").Value = 3264.2727
$ws.Range("K27").Value = 1771.1428
$ws.Range("L27").Value = 3264.2727
$ws.Range("M27").Value = -1664.1428
$ws.Range("N27").Value = -3478.2727

$ws.Range("H68").Value = 3150
$ws.Range("I68").Value = 3150
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3150
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2401
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 3150
$ws.Range("I71").Value = 3150
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 15750
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -12006
$ws.Range("N71").Value = ""

$ws.Range("H132").Value = 6250
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 80000
$ws.Range("I64").Value = 80000
$ws.Range("K64").Value = 80000
$ws.Range("M64").Value = -79752

$ws.Range("H67").Value = 80000
$ws.Range("I67").Value = 80000
$ws.Range("K67").Value = 80000
$ws.Range("M67").Value = -79142

$ws.Range("H107").Value = 3666.5
$ws.Range("J107").Value = 3666.5
$ws.Range("L107").Value = 10999.5
$ws.Range("N107").Value = -14839.5

$ws.Range("H122").Value = 7999
$ws.Range("I122").Value = 7999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 23997
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -21547
$ws.Range("N122").Value = ""

$ws.Range("H126").Value = 6750
$ws.Range("I126").Value = 5700
$ws.Range("K126").Value = 17100
$ws.Range("M126").Value = -14630
